# Quarterly cash-flow database update: drop the oldest quarter column (D,
# "فصل دوم منتهی به 1399/06") and append the newest quarter
# ("فصل چهارم منتهی به 1401/12" / published 1402-02-10) as the new last
# column (M). All the intervening quarter columns shift one position to
# the left, which EntireColumn.Delete gives us "for free" (cell values,
# shared-string refs and column widths all shift together, and orphaned
# shared strings are pruned automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest quarter (column D). Everything to the right slides
#    left by one column: E->D, F->E, ... M->L.
$ws.Range("D1").EntireColumn.Delete()

# 2) The freshly vacated column M needs the same custom width (31) as the
#    other "publish date" columns. ColumnWidth (characters) round-trips to
#    the stored `width` with a fixed +5/6 offset in this engine, so back
#    that out to land exactly on 31.
$ws.Range("M1").EntireColumn.ColumnWidth = 30.16666666666667

# 3) The old publish-date placeholder in J9 ("1401-10-28 (6)") is the same
#    shared string the post-shift I9 cell will end up pointing at, so edit
#    its text in place *before* the shift to land on the new value there.
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "1402-02-10 (7)"

# 4) New header labels for the newly-appended quarter. Force text format so
#    the ISO-looking date string isn't auto-parsed into a serial number.
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-10"

# 4) New financial figures for the appended quarter (column M).
$ws.Range("M12").Value = 7894154
$ws.Range("M13").Value = 78
$ws.Range("M14").Value = 7894232
$ws.Range("M17").Value = -221088
$ws.Range("M20").Value = -1500
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = -898865
$ws.Range("M25").Value = -2167529
$ws.Range("M26").Value = -5838436
$ws.Range("M31").Value = 1164967
$ws.Range("M32").Value = -7962451
$ws.Range("M33").Value = -68219
$ws.Range("M37").Value = 0
$ws.Range("M38").Value = 0
$ws.Range("M50").Value = -36091
$ws.Range("M51").Value = -36091
$ws.Range("M52").Value = -104310
$ws.Range("M53").Value = 2462928
$ws.Range("M54").Value = 11383
$ws.Range("M55").Value = 2370001
$ws.Range("M56").Value = 136334

# 5) Row 36 ("سود سهام پرداختی" offsets) breaks the pure left-shift
#    pattern: the algorithm now reports an actual 0 instead of the "-"
#    placeholder for both the I and M columns.
$ws.Range("I36").Value = 0
$ws.Range("M36").Value = 0

Write-Output "done"
